# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) for the affected Leve rows across multiple
# crafting-job sheets, per the latest data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3003
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 3003
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = $null
$ws.Range("M137").Value = 9009
$ws.Range("N137").Value = -14109
$ws.Range("H138").Value = 2596.3
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2367.5881
$ws.Range("I32").Value = 2367.5881
$ws.Range("K32").Value = 2367.5881
$ws.Range("M32").Value = -2080.5881
$ws.Range("H97").Value = 864.0833
$ws.Range("I97").Value = 864.0833
$ws.Range("K97").Value = 864.0833
$ws.Range("M97").Value = -368.0833
$ws.Range("H122").Value = 2422.3684
$ws.Range("I122").Value = 2273.9
$ws.Range("J122").Value = 2587.3333
$ws.Range("K122").Value = 6821.700000000001
$ws.Range("L122").Value = 7761.999899999999
$ws.Range("M122").Value = -4371.700000000001
$ws.Range("N122").Value = -12661.9999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1528.2106
$ws.Range("I86").Value = 1446.1333
$ws.Range("J86").Value = 1836
$ws.Range("K86").Value = 1446.1333
$ws.Range("L86").Value = 1836
$ws.Range("M86").Value = -323.1333
$ws.Range("N86").Value = -4082
$ws.Range("H89").Value = 1528.2106
$ws.Range("I89").Value = 1446.1333
$ws.Range("J89").Value = 1836
$ws.Range("K89").Value = 7230.666499999999
$ws.Range("L89").Value = 9180
$ws.Range("M89").Value = -1614.666499999999
$ws.Range("N89").Value = -20412
$ws.Range("H99").Value = 1338
$ws.Range("I99").Value = 1131.7142
$ws.Range("K99").Value = 1131.7142
$ws.Range("M99").Value = 366.2858000000001
$ws.Range("H105").Value = 3374.25
$ws.Range("I105").Value = 3374.25
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3374.25
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = $null
$ws.Range("N105").Value = -1627.25
$ws.Range("H107").Value = 5737.7
$ws.Range("I107").Value = 4280
$ws.Range("K107").Value = 4280
$ws.Range("M107").Value = -2360
$ws.Range("H134").Value = 3129.423
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2515.8948
$ws.Range("I16").Value = 1954.3846
$ws.Range("K16").Value = 1954.3846
$ws.Range("M16").Value = -1667.3846
$ws.Range("H31").Value = 2675.1667
$ws.Range("I31").Value = 5000
$ws.Range("J31").Value = 2210.2
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 2210.2
$ws.Range("M31").Value = -4705
$ws.Range("N31").Value = -2800.2
$ws.Range("H34").Value = 2675.1667
$ws.Range("I34").Value = 5000
$ws.Range("J34").Value = 2210.2
$ws.Range("K34").Value = 5000
$ws.Range("L34").Value = 2210.2
$ws.Range("M34").Value = -4798
$ws.Range("N34").Value = -2614.2
$ws.Range("H62").Value = 4747.5
$ws.Range("I62").Value = 4700
$ws.Range("J62").Value = 4795
$ws.Range("K62").Value = 4700
$ws.Range("L62").Value = 4795
$ws.Range("M62").Value = -4076
$ws.Range("N62").Value = -6043
$ws.Range("H65").Value = 4747.5
$ws.Range("I65").Value = 4700
$ws.Range("J65").Value = 4795
$ws.Range("K65").Value = 23500
$ws.Range("L65").Value = 23975
$ws.Range("M65").Value = -20380
$ws.Range("N65").Value = -30215
$ws.Range("H107").Value = 991.8
$ws.Range("I107").Value = 722.8182
$ws.Range("J107").Value = 1731.5
$ws.Range("K107").Value = 722.8182
$ws.Range("L107").Value = 1731.5
$ws.Range("M107").Value = 1197.1818
$ws.Range("N107").Value = -5571.5
$ws.Range("H113").Value = 2515.8948
$ws.Range("I113").Value = 1954.3846
$ws.Range("K113").Value = 1954.3846
$ws.Range("M113").Value = 215.6153999999999
$ws.Range("H132").Value = 1717.3928
$ws.Range("I132").Value = 1705.85
$ws.Range("J132").Value = 1746.25
$ws.Range("K132").Value = 5117.549999999999
$ws.Range("L132").Value = 5238.75
$ws.Range("M132").Value = -2587.549999999999
$ws.Range("N132").Value = -10298.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = $null
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = $null
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = 0
$ws.Range("H117").Value = 25237.25
$ws.Range("J117").Value = 50199.5
$ws.Range("L117").Value = 150598.5
$ws.Range("N117").Value = -157482.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 839.0909
$ws.Range("I16").Value = 903.2
$ws.Range("J16").Value = 198
$ws.Range("K16").Value = 903.2
$ws.Range("L16").Value = 198
$ws.Range("M16").Value = -733.2
$ws.Range("N16").Value = -538
$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -5272
$ws.Range("H93").Value = 1977.1428
$ws.Range("J93").Value = 3000
$ws.Range("L93").Value = 3000
$ws.Range("N93").Value = -5496
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 14287052
$ws.Range("I100").Value = 16668028
$ws.Range("K100").Value = 33336056
$ws.Range("M100").Value = -33335515
$ws.Range("H122").Value = 1224.75
$ws.Range("I122").Value = 824
$ws.Range("J122").Value = 2961.3333
$ws.Range("K122").Value = 2472
$ws.Range("L122").Value = 8883.999899999999
$ws.Range("M122").Value = -22
$ws.Range("N122").Value = -13783.9999
$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820
